# Applies the "second code review" edit to Data/Config.xlsx:
#  - Inserts a new row above the current row 29 on the "Settings" sheet
#    (pushes the existing DictionaryTransactionItemsColumns / SearchingForColumn /
#    DictionarySenderEmailColumn rows down from 29-32 to 31-33, with a blank
#    separator row left at 30 - matching the sheet's existing blank-row spacing
#    convention).
#  - Fills the newly freed row 29 with the new OutlookMailFolder setting.
#  - Moves the active selection to the newly added row (A29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a whole new row above row 29; everything from 29 downward shifts
# down by one (xlShiftDown), matching Excel's normal "Insert" row behaviour.
$ws.Rows("29:29").Insert()

# The freshly inserted row comes back at the default row height; restore the
# 14.25pt custom height used throughout the rest of this sheet.
$ws.Rows("29:29").RowHeight = 14.25

# Populate the newly inserted (now blank) row 29 with the new setting.
$ws.Range("A29").Value = "OutlookMailFolder"
$ws.Range("B29").Value = "Inbox"
$ws.Range("C29").Value = "It will look into this mail folder to read the mail messages."

# Update the sheet's saved selection/active cell to the new row, as recorded
# in the workbook after the edit.
$ws.Range("A29").Select() | Out-Null
